# daily auto push: 2026-03-01 03:22 UTC
# Insert a new data row for 2026/03/01 07:00 right after the existing
# 2026/03/01 04:00 entry (row 897), pushing the old row 898
# ("2026/12/29 13:00") and everything after it down by one row (old last
# row 939 becomes row 940).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 898 and everything below it down by one row.
$ws.Rows.Item(898).Insert()

# The date/weekday columns (A/B) are plain text in this sheet. Row 897
# already holds the same date/weekday text we need ("2026/03/01" / "日"),
# so copy it down as values instead of re-typing it - this avoids Excel
# auto-converting a typed date-like string into a real date serial and
# keeps the new cells free of any extra number-format/style.
$ws.Range("A897:B897").Copy()
$ws.Range("A898").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("C898").Value = 7
$ws.Range("D898").Value = 201
